$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("UNIVERSIDAD ESTATAL ADISTANCIA", $false, $false, $false, $false, $false, $true, 1, $false, "UNIVERSIDAD ESTATAL A DISTANCIA", 2)
Write-Output "found1=$found"
